# Rename sheets to remove diacritics/spaces (no-signs naming convention)
# "KPIs ngay" -> "KPI_ngay"
# "Tổng hợp" -> "Tong_hop"
# "TLXLTB và PH" -> "TLXLTB_PH"
# "Biểu đồ" -> "Bieu_Do"

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("KPIs ngay").Name = "KPI_ngay"
$wb.Worksheets.Item("Tổng hợp").Name = "Tong_hop"
$wb.Worksheets.Item("TLXLTB và PH").Name = "TLXLTB_PH"
$wb.Worksheets.Item("Biểu đồ").Name = "Bieu_Do"

# The chart on the renamed "Bieu_Do" sheet references that sheet by name in
# its series formulas; update them so the chart keeps pointing at the right
# ranges using the new (unquoted) sheet name.
$ws = $wb.Worksheets.Item("Bieu_Do")
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(`"SLPA`",Bieu_Do!`$C`$13:`$U`$13,Bieu_Do!`$C`$14:`$U`$14,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(`"TLPA/10K TB`",Bieu_Do!`$C`$13:`$U`$13,Bieu_Do!`$C`$15:`$U`$15,2)"
